# Auto-generated edit script: refreshes the live market-data-derived
# columns (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across all eight
# job sheets, matching the scheduled Sheets runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 461.53845
$ws.Range("I41").Value = 99
$ws.Range("J41").Value = 527.4545000000001
$ws.Range("K41").Value = 99
$ws.Range("L41").Value = 527.4545000000001
$ws.Range("M41").Value = 341
$ws.Range("N41").Value = -1407.4545
$ws.Range("H70").Value = 3460.3125
$ws.Range("H73").Value = 3460.3125
$ws.Range("H76").Value = 3155.889
$ws.Range("I76").Value = 3191.875
$ws.Range("J76").Value = 2868
$ws.Range("K76").Value = 3191.875
$ws.Range("L76").Value = 2868
$ws.Range("M76").Value = -2876.875
$ws.Range("N76").Value = -3498
$ws.Range("H79").Value = 3155.889
$ws.Range("I79").Value = 3191.875
$ws.Range("J79").Value = 2868
$ws.Range("K79").Value = 3191.875
$ws.Range("L79").Value = 2868
$ws.Range("M79").Value = -2099.875
$ws.Range("N79").Value = -5052
$ws.Range("H100").Value = 11766138
$ws.Range("I100").Value = 13334843
$ws.Range("J100").Value = 852
$ws.Range("K100").Value = 13334843
$ws.Range("L100").Value = 852
$ws.Range("M100").Value = -13334302
$ws.Range("N100").Value = -1934
$ws.Range("H103").Value = 17580
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 42750
$ws.Range("K103").Value = 2400
$ws.Range("L103").Value = 128250
$ws.Range("M103").Value = -1814
$ws.Range("N103").Value = -129422
$ws.Range("H112").Value = 438292.62
$ws.Range("I112").Value = 717.5
$ws.Range("J112").Value = 482050.12
$ws.Range("K112").Value = 2152.5
$ws.Range("L112").Value = 1446150.36
$ws.Range("M112").Value = -1044.5
$ws.Range("N112").Value = -1448366.36
$ws.Range("H132").Value = 436093.78
$ws.Range("I132").Value = 734.3158
$ws.Range("J132").Value = 2504051.2
$ws.Range("K132").Value = 2202.9474
$ws.Range("L132").Value = 7512153.600000001
$ws.Range("M132").Value = 327.0526
$ws.Range("N132").Value = -7517213.600000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3330.6711
$ws.Range("I32").Value = 2924.0447
$ws.Range("K32").Value = 2924.0447
$ws.Range("M32").Value = -2637.0447
$ws.Range("H61").Value = 1228.3914
$ws.Range("I61").Value = 1268.9143
$ws.Range("J61").Value = 1099.4546
$ws.Range("K61").Value = 1268.9143
$ws.Range("L61").Value = 1099.4546
$ws.Range("M61").Value = -1056.9143
$ws.Range("N61").Value = -1523.4546
$ws.Range("H63").Value = 7293048
$ws.Range("I63").Value = 10657041
$ws.Range("J63").Value = 4396.1665
$ws.Range("K63").Value = 10657041
$ws.Range("L63").Value = 4396.1665
$ws.Range("M63").Value = -10656355
$ws.Range("N63").Value = -5768.1665
$ws.Range("H66").Value = 7293048
$ws.Range("I66").Value = 10657041
$ws.Range("J66").Value = 4396.1665
$ws.Range("K66").Value = 53285205
$ws.Range("L66").Value = 21980.8325
$ws.Range("M66").Value = -53281773
$ws.Range("N66").Value = -28844.8325
$ws.Range("H74").Value = 3263.195
$ws.Range("I74").Value = 4206.864
$ws.Range("J74").Value = 2170.5264
$ws.Range("K74").Value = 4206.864
$ws.Range("L74").Value = 2170.5264
$ws.Range("M74").Value = -3332.864
$ws.Range("N74").Value = -3918.5264
$ws.Range("H77").Value = 3263.195
$ws.Range("I77").Value = 4206.864
$ws.Range("J77").Value = 2170.5264
$ws.Range("K77").Value = 21034.32
$ws.Range("L77").Value = 10852.632
$ws.Range("M77").Value = -16666.32
$ws.Range("N77").Value = -19588.632
$ws.Range("H136").Value = 1228.3914
$ws.Range("I136").Value = 1268.9143
$ws.Range("J136").Value = 1099.4546
$ws.Range("K136").Value = 3806.7429
$ws.Range("L136").Value = 3298.3638
$ws.Range("M136").Value = -1256.7429
$ws.Range("N136").Value = -8398.363799999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7042.5557
$ws.Range("I20").Value = 1365.4286
$ws.Range("J20").Value = 13156.385
$ws.Range("K20").Value = 1365.4286
$ws.Range("L20").Value = 13156.385
$ws.Range("M20").Value = -1118.4286
$ws.Range("N20").Value = -13650.385
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H105").Value = 1532.5
$ws.Range("I105").Value = 1529.4762
$ws.Range("J105").Value = 1596
$ws.Range("K105").Value = 1529.4762
$ws.Range("L105").Value = 1596
$ws.Range("M105").Value = 217.5237999999999
$ws.Range("N105").Value = -5090
$ws.Range("H107").Value = 1382
$ws.Range("I107").Value = 1382
$ws.Range("K107").Value = 1382
$ws.Range("M107").Value = 538
$ws.Range("H134").Value = 2490.375
$ws.Range("I134").Value = 1646.9
$ws.Range("K134").Value = 4940.700000000001
$ws.Range("M134").Value = -2405.700000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2855.44
$ws.Range("I31").Value = 1079.9166
$ws.Range("J31").Value = 4494.385
$ws.Range("K31").Value = 1079.9166
$ws.Range("L31").Value = 4494.385
$ws.Range("M31").Value = -784.9166
$ws.Range("N31").Value = -5084.385
$ws.Range("H34").Value = 2855.44
$ws.Range("I34").Value = 1079.9166
$ws.Range("J34").Value = 4494.385
$ws.Range("K34").Value = 1079.9166
$ws.Range("L34").Value = 4494.385
$ws.Range("M34").Value = -877.9166
$ws.Range("N34").Value = -4898.385
$ws.Range("H58").Value = 2077.4614
$ws.Range("I58").Value = 1853.9166
$ws.Range("J58").Value = 4760
$ws.Range("K58").Value = 1853.9166
$ws.Range("L58").Value = 4760
$ws.Range("M58").Value = -1650.9166
$ws.Range("N58").Value = -5166
$ws.Range("H122").Value = 2654.5
$ws.Range("I122").Value = 1080
$ws.Range("J122").Value = 3779.1428
$ws.Range("K122").Value = 3240
$ws.Range("L122").Value = 11337.4284
$ws.Range("M122").Value = -790
$ws.Range("N122").Value = -16237.4284
$ws.Range("H136").Value = 2077.4614
$ws.Range("I136").Value = 1853.9166
$ws.Range("J136").Value = 4760
$ws.Range("K136").Value = 5561.7498
$ws.Range("L136").Value = 14280
$ws.Range("M136").Value = -3011.7498
$ws.Range("N136").Value = -19380

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 727.55554
$ws.Range("I113").Value = 720.48
$ws.Range("K113").Value = 2161.44
$ws.Range("M113").Value = 8.559999999999945
$ws.Range("H137").Value = 1394
$ws.Range("I137").Value = 992.5
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 2977.5
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = 2122.5
$ws.Range("N137").Value = -19200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6214.969
$ws.Range("I70").Value = 5805.9624
$ws.Range("J70").Value = 8021.4165
$ws.Range("K70").Value = 5805.9624
$ws.Range("L70").Value = 8021.4165
$ws.Range("M70").Value = -5535.9624
$ws.Range("N70").Value = -8561.416499999999
$ws.Range("H73").Value = 6214.969
$ws.Range("I73").Value = 5805.9624
$ws.Range("J73").Value = 8021.4165
$ws.Range("K73").Value = 5805.9624
$ws.Range("L73").Value = 8021.4165
$ws.Range("M73").Value = -4869.9624
$ws.Range("N73").Value = -9893.416499999999
$ws.Range("H132").Value = 2262.525
$ws.Range("I132").Value = 1695.8182
$ws.Range("J132").Value = 2955.1667
$ws.Range("K132").Value = 5087.4546
$ws.Range("L132").Value = 8865.500100000001
$ws.Range("M132").Value = -2557.4546
$ws.Range("N132").Value = -13925.5001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 29985
$ws.Range("J26").Value = 29985
$ws.Range("L26").Value = 29985
$ws.Range("N26").Value = -30575
$ws.Range("H93").Value = 4117619.2
$ws.Range("I93").Value = 7409785.5
$ws.Range("J93").Value = 2411.5
$ws.Range("K93").Value = 7409785.5
$ws.Range("L93").Value = 2411.5
$ws.Range("M93").Value = -7408537.5
$ws.Range("N93").Value = -4907.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12822444
$ws.Range("I132").Value = 1110.95
$ws.Range("K132").Value = 3332.85
$ws.Range("M132").Value = -802.8500000000004
